# Corrigindo erro de Lógica NRU
# Fixes the cached/entered NRU (Não Usado Recentemente) results in the
# "Planilha1" tables so they match the corrected simulation logic, and
# marks the end of the last data block with a thin left/right border.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- Table 1: "11 Páginas 3 Frames" (rows 3-12) ---------------------------
# NRU column (T) was 11 for every run; corrected to 9.
$ws.Range("T3:T12").Value = 9

# --- Table 2: "10 Páginas 3 Frames" (rows 16-26 / 16-25) -------------------
# NRU column (T) was 10 for every run; corrected to 7 (rows 17-26).
$ws.Range("T17:T26").Value = 7
# Second table's columns (Z = second_chance, AA = NRU) in the same block.
$ws.Range("Z16:Z25").Value = 9
$ws.Range("AA16:AA25").Value = 7

# --- Table 3: "10 Páginas 2 Frames" (rows 30-39) ---------------------------
# NRU column (T) was 7 for every run; corrected to 6.
$ws.Range("T30:T39").Value = 6
# Second table's NRU column (AA) in the same block.
$ws.Range("AA30:AA39").Value = 6

# --- New separator row below the last table --------------------------------
# Row 40 / cell T40: empty cell carrying a thin left+right border.
$sep = $ws.Range("T40")
$sep.Value = ""
$sep.Borders.Item(7).LineStyle = 1
$sep.Borders.Item(7).Weight = 2
$sep.Borders.Item(10).LineStyle = 1
$sep.Borders.Item(10).Weight = 2

# --- Selection / view bookkeeping ------------------------------------------
$ws.Activate()
$ws.Range("AE10").Select()
